# Update DateBase/orders/Fresh bloom Flowers_2025-10-15.xlsx
# - Append rows 42-51 to the "Orders" sheet (columns A/C/F)
# - Update the summary "Number" digest string in G2 of the "Summary" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

function Set-TextValue {
    param($Sheet, $Addr, $Text)
    $rng = $Sheet.Range($Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
}

$newRows = @(
    @{ Row = 41; A = $null; C = $null; F = "5" },
    @{ Row = 42; A = $null; C = "321_雪柳叶_Spiraea  leaves_undefined_1bunch"; F = "40" },
    @{ Row = 43; A = $null; C = "439_九星叶_undefined_undefined_1bunch"; F = "12" },
    @{ Row = 44; A = $null; C = "300_白星_White Gypso_ gypsophila_1kg"; F = "4" },
    @{ Row = 45; A = $null; C = "302_彩星 浅粉_Tinted Gypso light pink_undefined_0.5kg"; F = "20" },
    @{ Row = 46; A = "17"; C = "349_千层金绿_Melaleuca bracteata`n（dyed orange）_Melaleuca bracteata F.Muell._1bunch"; F = "10" },
    @{ Row = 47; A = $null; C = "401_大飞燕白色_delphinium white_undefined_1bunch"; F = "40" },
    @{ Row = 48; A = $null; C = "403_大飞燕浅蓝色_delphinium light blue_undefined_1bunch"; F = "25" },
    @{ Row = 49; A = $null; C = "107_绣球单瓣浅粉_Hydrangea Light Pink S_Hydrangea L._1stem"; F = "15" },
    @{ Row = 50; A = $null; C = "100_绣球单瓣白_Hydrangea White S_Hydrangea L._1stem"; F = "60" },
    @{ Row = 51; A = $null; C = "118_绣球老绿_Hydrangea Garden Lace_Hydrangea L._1stem"; F = $null }
)

foreach ($r in $newRows) {
    if ($r.A) {
        $addrA = "A" + $r.Row
        Set-TextValue $ws $addrA $r.A
    }
    if ($r.C) {
        $addrC = "C" + $r.Row
        Set-TextValue $ws $addrC $r.C
    }
    if ($r.F) {
        $addrF = "F" + $r.Row
        Set-TextValue $ws $addrF $r.F
    }
}

$ws2 = $wb.Worksheets.Item("Summary")
$g2 = "015196181942320232115225241410308117766324040401156054208556562355512205401242010402515600"
Set-TextValue $ws2 "G2" $g2
